$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "usd_volume" placeholder header cell in the template row
$ws.Range("M2").Value = '${row.usd_volume}'

# Set the column M width to match the new narrower content (was 20.33 -> 8.82)
$ws.Columns.Item(13).ColumnWidth = 8.0

# Update the view's top-left visible cell and active selection
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("M3").Select()
